$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that wraps the second scatter-plot image
#    (the one in the "Repeat #2 above for attribute2 and attribute4" answer).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) Locate the "Min leaf node size = 5" paragraph and the blank paragraph
#    right after it (which currently has ind left=1545 / 77.25pt).
# ---------------------------------------------------------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Min leaf node size = 5") {
        $targetIndex = $i
        break
    }
}

$blankAfter5 = $d.Paragraphs.Item($targetIndex + 1)
# change ind left from 1545 twips (77.25pt) to -1170 twips (-58.5pt)
$blankAfter5.LeftIndent = -58.5

# ---------------------------------------------------------------------------
# 3) Insert the new paragraphs right after that blank paragraph:
#    - "Min leaf node size = 25" (ilvl=1, numId=1)
#    - blank (ListParagraph, no ind)
#    - blank (ListParagraph, ind left=825 / 41.25pt)
#    - blank (ListParagraph, no ind)
#    - "Min leaf node size = 50" (ilvl=1, numId=1)
#    - a totally empty paragraph
# ---------------------------------------------------------------------------
$anchor = $blankAfter5

$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($targetIndex + 2)
$p1.Style = "List Paragraph"
$p1.Range.ListFormat.ListLevelNumber = 2
$p1.Range.Text = "Min leaf node size = 25"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($targetIndex + 3)
$p2.Style = "List Paragraph"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($targetIndex + 4)
$p3.Style = "List Paragraph"
$p3.LeftIndent = 41.25

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($targetIndex + 5)
$p4.Style = "List Paragraph"

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item($targetIndex + 6)
$p5.Style = "List Paragraph"
$p5.Range.ListFormat.ListLevelNumber = 2
$p5.Range.Text = "Min leaf node size = 50"

$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item($targetIndex + 7)
$p6.Style = "Normal"

# ---------------------------------------------------------------------------
# 4) Split the final "(40) Test each..." paragraph's text and re-insert the
#    _GoBack bookmark at the split point (Word re-adds this automatically
#    whenever the cursor/edit was last at that spot).
# ---------------------------------------------------------------------------
$finalIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("(40) Test each")) {
        $finalIndex = $i
        break
    }
}
$finalPara = $d.Paragraphs.Item($finalIndex)
$splitPoint = $finalPara.Range.Start + 300
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange)

Write-Output "done"
